# Insert a new record at row 98 of the "Hortaliza, Macroferia Regional de
# Talca - Zapallo" data table. Existing rows 98..202 shift down to 99..203
# (EntireRow insert handles that, carrying their data/formatting with them),
# and the new row 98 is populated with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data down by one row, starting at row 98.
$ws.Rows(98).Insert()

# Fill in the newly inserted row 98 with the new record's values.
$ws.Range("A98").Value = 5
$ws.Range("B98").Value = "Macroferia Regional de Talca"
$ws.Range("C98").Value = "Maule"
$ws.Range("D98").Value = 44629
$ws.Range("E98").Value = 7
$ws.Range("F98").Value = 100112045
$ws.Range("G98").Value = "Zapallo"
$ws.Range("H98").Value = "Camote"
$ws.Range("I98").Value = "1a nueva(o)"
$ws.Range("J98").Value = 900
$ws.Range("K98").Value = 300
$ws.Range("L98").Value = 300
$ws.Range("M98").Value = 300
$ws.Range("N98").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O98").Value = "Región del Maule"
$ws.Range("P98").Value = 300
$ws.Range("Q98").Value = 1
$ws.Range("R98").Value = "Hortaliza"
